$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on D-column cells whose new numeric-looking values
# (e.g. "1.003", "315.50") would otherwise be auto-converted to numbers
# by Excel, corrupting trailing zeros / precision / significant digits.
$textCells = @("D4","D5","D7","D9","D10","D11","D12","D13","D14","D15","D16","D18","D19","D21","D23","D25","D26","D27","D28","D29","D30","D31","D32","D34","D35","D36","D37","D38","D39","D40","D41","D42","D43","D44","D45","D46","D47","D48","D49","D50","D51")
foreach ($cellRef in $textCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$ws.Range('D2').Value = '24.916.11'
$ws.Range('E2').Value = '  +0.42%  '

$ws.Range('D3').Value = '1.692.25'
$ws.Range('E3').Value = '  -1.02%  '

$ws.Range('D4').Value = '1.003'
$ws.Range('E4').Value = '  +1.16%  '

$ws.Range('D5').Value = '315.50'
$ws.Range('E5').Value = '  -0.42%  '

$ws.Range('E6').Value = '  +0.94%  '

$ws.Range('D7').Value = '0.3945'
$ws.Range('E7').Value = '  +0.84%  '

$ws.Range('E8').Value = '  -1.89%  '

$ws.Range('B9').Value = 'OKB'
$ws.Range('C9').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D9').Value = '52.83'
$ws.Range('E9').Value = '  -2.50%  '

$ws.Range('B10').Value = 'Polygon'
$ws.Range('C10').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D10').Value = '1.443'
$ws.Range('E10').Value = '  -3.48%  '

$ws.Range('D11').Value = '1.005'
$ws.Range('E11').Value = '  +0.74%  '

$ws.Range('D12').Value = '0.08721'
$ws.Range('E12').Value = '  -1.31%  '

$ws.Range('D13').Value = '25.46'
$ws.Range('E13').Value = '  -2.92%  '

$ws.Range('D14').Value = '7.375'
$ws.Range('E14').Value = '  -1.53%  '

$ws.Range('B15').Value = 'Chainlink'
$ws.Range('C15').Value = 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
$ws.Range('D15').Value = '7.867'
$ws.Range('E15').Value = '  -3.60%  '

$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').Value = '0.00001336'
$ws.Range('E16').Value = '  -2.05%  '

$ws.Range('D17').Value = '1.742.42'
$ws.Range('E17').Value = '  +2.41%  '

$ws.Range('D18').Value = '94.81'
$ws.Range('E18').Value = '  -3.27%  '

$ws.Range('D19').Value = '0.07190'
$ws.Range('E19').Value = '  +0.16%  '

$ws.Range('E20').Value = '  -1.30%  '

$ws.Range('D21').Value = '7.173'
$ws.Range('E21').Value = '  -2.10%  '

$ws.Range('E22').Value = '  +0.93%  '

$ws.Range('D23').Value = '14.17'
$ws.Range('E23').Value = '  -1.49%  '

$ws.Range('D24').Value = '24.879.92'
$ws.Range('E24').Value = '  +0.34%  '

$ws.Range('D25').Value = '2.408'
$ws.Range('E25').Value = '  +3.33%  '

$ws.Range('D26').Value = '2.844'
$ws.Range('E26').Value = '  -6.43%  '

$ws.Range('D27').Value = '23.04'
$ws.Range('E27').Value = '  -0.28%  '

$ws.Range('D28').Value = '6.036'
$ws.Range('E28').Value = '  +1.60%  '

$ws.Range('D29').Value = '162.24'
$ws.Range('E29').Value = '  -3.06%  '

$ws.Range('D30').Value = '148.58'
$ws.Range('E30').Value = '  +2.98%  '

$ws.Range('D31').Value = '2.658'
$ws.Range('E31').Value = '  +22.12%  '

$ws.Range('D32').Value = '8.049'
$ws.Range('E32').Value = '  -5.25%  '

$ws.Range('D33').Value = '1.898.88'
$ws.Range('E33').Value = '  +0.58%  '

$ws.Range('D34').Value = '0.08521'
$ws.Range('E34').Value = '  -3.51%  '

$ws.Range('D35').Value = '0.03110'
$ws.Range('E35').Value = '  -0.07%  '

$ws.Range('D36').Value = '1.033'
$ws.Range('E36').Value = '  -2.89%  '

$ws.Range('B37').Value = 'Algorand'
$ws.Range('C37').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D37').Value = '0.2862'
$ws.Range('E37').Value = '  +1.58%  '

$ws.Range('B38').Value = 'InternetComputer(DFINITY)'
$ws.Range('C38').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D38').Value = '7.001'
$ws.Range('E38').Value = '  -3.99%  '

$ws.Range('D39').Value = '0.09687'
$ws.Range('E39').Value = '  +5.24%  '

$ws.Range('D40').Value = '10.76'
$ws.Range('E40').Value = '  -1.84%  '

$ws.Range('D41').Value = '0.8065'
$ws.Range('E41').Value = '  -8.43%  '

$ws.Range('D42').Value = '13.84'
$ws.Range('E42').Value = '  -3.23%  '

$ws.Range('D43').Value = '1.468'
$ws.Range('E43').Value = '  -1.15%  '

$ws.Range('D44').Value = '16.91'
$ws.Range('E44').Value = '  -2.91%  '

$ws.Range('D45').Value = '2.624'
$ws.Range('E45').Value = '  -2.42%  '

$ws.Range('D46').Value = '0.7270'
$ws.Range('E46').Value = '  -4.01%  '

$ws.Range('D47').Value = '4.215'
$ws.Range('E47').Value = '  -0.92%  '

$ws.Range('D48').Value = '0.08935'
$ws.Range('E48').Value = '  +7.93%  '

$ws.Range('D49').Value = '1.387'
$ws.Range('E49').Value = '  -1.24%  '

$ws.Range('D50').Value = '1.003'
$ws.Range('E50').Value = '  +0.82%  '

$ws.Range('D51').Value = '139.29'
$ws.Range('E51').Value = '  -1.22%  '
